# "test P7 with -10 percent"
# Update the cached solver results across several sheets to reflect a
# re-run of the scenario (P7) with the capacity/parameter turned down by
# 10 percent. All target cells already exist; we simply overwrite their
# values with the new, recomputed results.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "general": summary KPIs (objective value, runtime, Z1)
# ---------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Range("B3").Value = 33.0671278619396
$wsGeneral.Range("B4").Value = 0.009000062942504883
$wsGeneral.Range("B6").Value = 33.0671278619396

# ---------------------------------------------------------------------
# Sheet "x": assignment permutation column (j per i)
# ---------------------------------------------------------------------
$wsX = $wb.Worksheets.Item("x")
$wsX.Range("B2").Value = 1
$wsX.Range("B3").Value = 3
$wsX.Range("B4").Value = 4
$wsX.Range("B6").Value = 12
$wsX.Range("B7").Value = 6
$wsX.Range("B8").Value = 8
$wsX.Range("B9").Value = 13
$wsX.Range("B10").Value = 2
$wsX.Range("B12").Value = 9
$wsX.Range("B13").Value = 5
$wsX.Range("B14").Value = 11

# ---------------------------------------------------------------------
# Sheet "TBar": per-stop time values
# ---------------------------------------------------------------------
$wsTBar = $wb.Worksheets.Item("TBar")
$wsTBar.Range("B3").Value = 30
$wsTBar.Range("B4").Value = 34.69770569366315
$wsTBar.Range("B5").Value = 30
$wsTBar.Range("B6").Value = 30
$wsTBar.Range("B7").Value = 36.71579249669672
$wsTBar.Range("B8").Value = 30.34885527085025
$wsTBar.Range("B9").Value = 30.60033324079214
$wsTBar.Range("B10").Value = 32.31224998648503
$wsTBar.Range("B11").Value = 34.76592070603971
$wsTBar.Range("B12").Value = 32.61192465059682
$wsTBar.Range("B13").Value = 36.71671453559702
$wsTBar.Range("B14").Value = 37.27819014430416
$wsTBar.Range("B15").Value = 37.90090852477161

# ---------------------------------------------------------------------
# Sheet "Q": per (j,s) flow quantities
# ---------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Range("C7").Value = 109.9450000000008
$wsQ.Range("C8").Value = 117.5900000000008
$wsQ.Range("C9").Value = 113.2700000000008
$wsQ.Range("C10").Value = 119.1550000000008
$wsQ.Range("C11").Value = 115.8050000000008
$wsQ.Range("C12").Value = 235.775
$wsQ.Range("C13").Value = 229.025
$wsQ.Range("C14").Value = 213.42
$wsQ.Range("C15").Value = 226.76
$wsQ.Range("C16").Value = 221.56
$wsQ.Range("C17").Value = 46.91999999999942
$wsQ.Range("C18").Value = 36.10499999999942
$wsQ.Range("C19").Value = 34.91499999999942
$wsQ.Range("C20").Value = 37.48999999999942
$wsQ.Range("C21").Value = 39.43499999999941
$wsQ.Range("C22").Value = 72.6299999999995
$wsQ.Range("C23").Value = 80.0549999999995
$wsQ.Range("C24").Value = 82.31999999999948
$wsQ.Range("C25").Value = 83.9549999999995
$wsQ.Range("C26").Value = 80.8149999999995
$wsQ.Range("C27").Value = 295.9199999999996
$wsQ.Range("C28").Value = 323.5
$wsQ.Range("C29").Value = 294.2649999999996
$wsQ.Range("C30").Value = 311.1
$wsQ.Range("C31").Value = 297.3649999999996
$wsQ.Range("C32").Value = 154.3
$wsQ.Range("C33").Value = 148.3449999999993
$wsQ.Range("C34").Value = 128.7049999999993
$wsQ.Range("C35").Value = 146.3249999999993
$wsQ.Range("C36").Value = 134.2149999999993
$wsQ.Range("C37").Value = 193.0200000000017
$wsQ.Range("C38").Value = 202.3100000000017
$wsQ.Range("C39").Value = 191.2450000000017
$wsQ.Range("C40").Value = 208.9250000000017
$wsQ.Range("C41").Value = 197.6600000000017
$wsQ.Range("C42").Value = 140.5549999999989
$wsQ.Range("C43").Value = 159.2149999999988
$wsQ.Range("C44").Value = 142.1399999999989
$wsQ.Range("C45").Value = 147.7249999999989
$wsQ.Range("C46").Value = 139.7449999999989
$wsQ.Range("C47").Value = 226.0399999999994
$wsQ.Range("C48").Value = 247.1799999999993
$wsQ.Range("C49").Value = 221.8549999999994
$wsQ.Range("C50").Value = 238.4549999999994
$wsQ.Range("C51").Value = 224.4749999999994
$wsQ.Range("C52").Value = 250.970000000001
$wsQ.Range("C53").Value = 260.9900000000009
$wsQ.Range("C54").Value = 252.975000000001
$wsQ.Range("C55").Value = 269.580000000001
$wsQ.Range("C56").Value = 250.575000000001
$wsQ.Range("C57").Value = 250.970000000001
$wsQ.Range("C58").Value = 260.9900000000009
$wsQ.Range("C59").Value = 252.975000000001
$wsQ.Range("C60").Value = 269.580000000001
$wsQ.Range("C61").Value = 250.575000000001
$wsQ.Range("C62").Value = 235.775
$wsQ.Range("C63").Value = 229.025
$wsQ.Range("C64").Value = 213.42
$wsQ.Range("C65").Value = 226.76
$wsQ.Range("C66").Value = 221.56
$wsQ.Range("C67").Value = 295.9199999999996
$wsQ.Range("C68").Value = 323.5
$wsQ.Range("C69").Value = 294.2649999999996
$wsQ.Range("C70").Value = 311.1
$wsQ.Range("C71").Value = 297.3649999999996
